# Auto-generated edit script
# Applies updated FFXIV market-price snapshot values (H,I,J,K,L,M,N columns)
# to the rows identified in the commit diff, per worksheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 14014742
$ws.Range("I106").Value = 16016561
$ws.Range("K106").Value = 16016561
$ws.Range("M106").Value = -16015930
$ws.Range("H129").Value = 1009.18335
$ws.Range("J129").Value = 1043.6072
$ws.Range("L129").Value = 3130.8216
$ws.Range("N129").Value = -13130.8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3148.434
$ws.Range("I32").Value = 2301.4565
$ws.Range("J32").Value = 8714.286
$ws.Range("K32").Value = 2301.4565
$ws.Range("L32").Value = 8714.286
$ws.Range("M32").Value = -2014.4565
$ws.Range("N32").Value = -9288.286
$ws.Range("H61").Value = 1845.7354
$ws.Range("I61").Value = 882.8261
$ws.Range("J61").Value = 3859.0908
$ws.Range("K61").Value = 882.8261
$ws.Range("L61").Value = 3859.0908
$ws.Range("M61").Value = -670.8261
$ws.Range("N61").Value = -4283.0908
$ws.Range("H74").Value = 1340.7878
$ws.Range("I74").Value = 1358.909
$ws.Range("K74").Value = 1358.909
$ws.Range("M74").Value = -484.9090000000001
$ws.Range("H77").Value = 1340.7878
$ws.Range("I77").Value = 1358.909
$ws.Range("K77").Value = 6794.545
$ws.Range("M77").Value = -2426.545
$ws.Range("H102").Value = 1413.6923
$ws.Range("I102").Value = 1120
$ws.Range("K102").Value = 1120
$ws.Range("M102").Value = 502
$ws.Range("H110").Value = 802.913
$ws.Range("I110").Value = 773
$ws.Range("J110").Value = 945
$ws.Range("K110").Value = 773
$ws.Range("L110").Value = 945
$ws.Range("M110").Value = 1272
$ws.Range("N110").Value = -5035
$ws.Range("H122").Value = 1854.4
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 2092
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 6276
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -11176
$ws.Range("H132").Value = 1753.3438
$ws.Range("I132").Value = 1390.6072
$ws.Range("J132").Value = 4292.5
$ws.Range("K132").Value = 4171.821599999999
$ws.Range("L132").Value = 12877.5
$ws.Range("M132").Value = -1641.821599999999
$ws.Range("N132").Value = -17937.5
$ws.Range("H136").Value = 1845.7354
$ws.Range("I136").Value = 882.8261
$ws.Range("J136").Value = 3859.0908
$ws.Range("K136").Value = 2648.4783
$ws.Range("L136").Value = 11577.2724
$ws.Range("M136").Value = -98.47829999999976
$ws.Range("N136").Value = -16677.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9971.154
$ws.Range("I86").Value = 1916.1428
$ws.Range("J86").Value = 19368.666
$ws.Range("K86").Value = 1916.1428
$ws.Range("L86").Value = 19368.666
$ws.Range("M86").Value = -793.1428000000001
$ws.Range("N86").Value = -21614.666
$ws.Range("H89").Value = 9971.154
$ws.Range("I89").Value = 1916.1428
$ws.Range("J89").Value = 19368.666
$ws.Range("K89").Value = 9580.714
$ws.Range("L89").Value = 96843.33
$ws.Range("M89").Value = -3964.714
$ws.Range("N89").Value = -108075.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1677.8572
$ws.Range("I31").Value = 1056.1875
$ws.Range("J31").Value = 2506.75
$ws.Range("K31").Value = 1056.1875
$ws.Range("L31").Value = 2506.75
$ws.Range("M31").Value = -761.1875
$ws.Range("N31").Value = -3096.75
$ws.Range("H34").Value = 1677.8572
$ws.Range("I34").Value = 1056.1875
$ws.Range("J34").Value = 2506.75
$ws.Range("K34").Value = 1056.1875
$ws.Range("L34").Value = 2506.75
$ws.Range("M34").Value = -854.1875
$ws.Range("N34").Value = -2910.75
$ws.Range("H58").Value = 2717.1667
$ws.Range("I58").Value = 969.5
$ws.Range("K58").Value = 969.5
$ws.Range("M58").Value = -766.5
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -10300
$ws.Range("H132").Value = 4405
$ws.Range("I132").Value = 3447.2
$ws.Range("K132").Value = 10341.6
$ws.Range("M132").Value = -7811.599999999999
$ws.Range("H134").Value = 3645.7778
$ws.Range("I134").Value = 1560.6
$ws.Range("K134").Value = 4681.799999999999
$ws.Range("M134").Value = -2146.799999999999
$ws.Range("H136").Value = 2717.1667
$ws.Range("I136").Value = 969.5
$ws.Range("K136").Value = 2908.5
$ws.Range("M136").Value = -358.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1028.5354
$ws.Range("I68").Value = 794.8393
$ws.Range("J68").Value = 1332.8837
$ws.Range("K68").Value = 2384.5179
$ws.Range("L68").Value = 3998.6511
$ws.Range("M68").Value = -1573.5179
$ws.Range("N68").Value = -5620.6511
$ws.Range("H71").Value = 1028.5354
$ws.Range("I71").Value = 794.8393
$ws.Range("J71").Value = 1332.8837
$ws.Range("K71").Value = 7153.553699999999
$ws.Range("L71").Value = 11995.9533
$ws.Range("M71").Value = -3097.553699999999
$ws.Range("N71").Value = -20107.9533
$ws.Range("H107").Value = 1081.6721
$ws.Range("I107").Value = 354.73685
$ws.Range("J107").Value = 1410.5238
$ws.Range("K107").Value = 1064.21055
$ws.Range("L107").Value = 4231.5714
$ws.Range("M107").Value = 855.78945
$ws.Range("N107").Value = -8071.5714
$ws.Range("H113").Value = 41667092
$ws.Range("I113").Value = 407.5
$ws.Range("J113").Value = 83333780
$ws.Range("K113").Value = 1222.5
$ws.Range("L113").Value = 250001340
$ws.Range("M113").Value = 947.5
$ws.Range("N113").Value = -250005680
$ws.Range("H122").Value = 918.4091
$ws.Range("I122").Value = 745.3
$ws.Range("J122").Value = 1062.6666
$ws.Range("K122").Value = 6707.7
$ws.Range("L122").Value = 9563.999400000001
$ws.Range("M122").Value = -4257.7
$ws.Range("N122").Value = -14463.9994
$ws.Range("H137").Value = 5946938.5
$ws.Range("J137").Value = 116837.336
$ws.Range("L137").Value = 350512.008
$ws.Range("N137").Value = -360712.008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2699.75
$ws.Range("I80").Value = 2650
$ws.Range("J80").Value = 2749.5
$ws.Range("K80").Value = 2650
$ws.Range("L80").Value = 2749.5
$ws.Range("M80").Value = -1652
$ws.Range("N80").Value = -4745.5
$ws.Range("H83").Value = 2699.75
$ws.Range("I83").Value = 2650
$ws.Range("J83").Value = 2749.5
$ws.Range("K83").Value = 13250
$ws.Range("L83").Value = 13747.5
$ws.Range("M83").Value = -8258
$ws.Range("N83").Value = -23731.5
$ws.Range("H113").Value = 1458.45
$ws.Range("I113").Value = 1249.5333
$ws.Range("J113").Value = 2085.2
$ws.Range("K113").Value = 1249.5333
$ws.Range("L113").Value = 2085.2
$ws.Range("M113").Value = 920.4666999999999
$ws.Range("N113").Value = -6425.2
$ws.Range("H126").Value = 1991.4193
$ws.Range("I126").Value = 1353.8182
$ws.Range("K126").Value = 4061.4546
$ws.Range("M126").Value = -1591.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 50000
$ws.Range("I4").Value = 50000
$ws.Range("K4").Value = 50000
$ws.Range("M4").Value = -49887
$ws.Range("H28").Value = 50000
$ws.Range("I28").Value = 50000
$ws.Range("K28").Value = 50000
$ws.Range("M28").Value = -49768
$ws.Range("H37").Value = 50000
$ws.Range("I37").Value = 50000
$ws.Range("K37").Value = 50000
$ws.Range("M37").Value = -49893
$ws.Range("H40").Value = 4459.375
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 4670
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 4670
$ws.Range("M40").Value = -1164
$ws.Range("N40").Value = -4942
$ws.Range("H122").Value = 3308.077
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 3500.4167
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 10501.2501
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -15401.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 38465556
$ws.Range("I132").Value = 62503800
$ws.Range("J132").Value = 4372.6
$ws.Range("K132").Value = 187511400
$ws.Range("L132").Value = 13117.8
$ws.Range("M132").Value = -187508870
$ws.Range("N132").Value = -18177.8
$ws.Range("H136").Value = 18520364
$ws.Range("I136").Value = 30303926
$ws.Range("J136").Value = 3335.2856
$ws.Range("K136").Value = 90911778
$ws.Range("L136").Value = 10005.8568
$ws.Range("M136").Value = -90909228
$ws.Range("N136").Value = -15105.8568
